$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.073.91"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.851.56"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "706.39"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.01"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.847.90"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.30"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.86"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.501.12"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.869.86"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.093.52"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.21"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.37"
$ws.Range("E21").Value = "  +3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.65"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.21"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.21"
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.53"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.50"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -4.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.19"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.808.36"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  +7.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.04"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  +5.82%  "
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000324"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.43"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.60"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "415.82"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("E51").Value = "  +0.36%  "
